$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the shared string "failure" -> "fail" (this affects cell F6 which uses that string)
$ws.Range("F6").Value = "fail"

# Row 2: F2, G2 from 1 -> 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0

# Row 3: add G3 = 0
$ws.Range("G3").Value = 0

# Row 4: add F4, G4, H4, I4
$ws.Range("F4").Value = 3.125
$ws.Range("G4").Value = 3.125
$ws.Range("H4").Value = 0.010631716963287633
$ws.Range("I4").Value = 0.0056566983573779128

# Row 5: F5, G5 from 1 -> 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0

# Row 8: F8, G8 from 1 -> 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0

# Row 9: F9, G9 from 1 -> 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0

# Update the selected cell/active cell to E8
$ws.Range("E8").Select()
